# Fix bug for ValueError in Add Employee function
# - A6 ("employee_id" for Thim) was stored as text "5"; correct it to the
#   number 5.
# - Append a new employee row (row 7): Xi, id 6, rate/hours/overtime all 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing employee_id that was mistakenly written as a string.
$ws.Cells.Item(6, 1).Value = 5

# Add the new employee row. employee_id stays text here (same pre-fix shape
# as the other rows), the numeric fields are real numbers.
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "6"
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).Value = "Xi"
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).Value = 10
